$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row 6 to the ledger, following the same pattern as the existing rows.
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = "Cash To Mustafa"
$ws.Range("D6").Value = 77000
$ws.Range("E6").Value = "MZN"
$ws.Range("F6").Value = 2
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 156777
$ws.Range("I6").Value = "USD"
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = "Cash To Mustafa"
$ws.Range("L6").Value = "money transfer to xyz"
